$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D (Price) cells as Text so numeric-looking
# strings (e.g. "203.36") are preserved exactly as typed instead of
# being coerced into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.538.83'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '3.527.67'
$ws.Range("E3").Value = '  -2.58%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '203.36'
$ws.Range("E5").Value = '  +4.11%  '
$ws.Range("D6").Value = '551.64'
$ws.Range("E6").Value = '  -6.41%  '
$ws.Range("D7").Value = '3.521.45'
$ws.Range("E7").Value = '  -2.57%  '
$ws.Range("E8").Value = '  -2.57%  '
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '0.657'
$ws.Range("E10").Value = '  -4.00%  '
$ws.Range("D11").Value = '62.57'
$ws.Range("E11").Value = '  +11.77%  '
$ws.Range("D12").Value = '0.142'
$ws.Range("E12").Value = '  -7.17%  '
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").Value = '  -7.37%  '
$ws.Range("E14").Value = '  -1.79%  '
$ws.Range("D15").Value = '4.082.64'
$ws.Range("E15").Value = '  -2.92%  '
$ws.Range("D16").Value = '3.525.44'
$ws.Range("E16").Value = '  -2.68%  '
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").Value = '18.58'
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '67.200.12'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("D20").Value = '11.90'
$ws.Range("E20").Value = '  -5.58%  '
$ws.Range("D21").Value = '1.03'
$ws.Range("E21").Value = '  -4.46%  '
$ws.Range("D22").Value = '391.08'
$ws.Range("E22").Value = '  -3.62%  '
$ws.Range("D23").Value = '4.03'
$ws.Range("E23").Value = '  -6.05%  '
$ws.Range("D24").Value = '11.98'
$ws.Range("E24").Value = '  -10.68%  '
$ws.Range("D25").Value = '82.71'
$ws.Range("E25").Value = '  -4.11%  '
$ws.Range("E26").Value = '  -4.96%  '
$ws.Range("D27").Value = '12.14'
$ws.Range("E27").Value = '  -3.64%  '
$ws.Range("D28").Value = '3.75'
$ws.Range("E28").Value = '  -4.07%  '
$ws.Range("D29").Value = '8.92'
$ws.Range("E29").Value = '  -3.65%  '
$ws.Range("D30").Value = '30.86'
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("D31").Value = '7.30'
$ws.Range("E31").Value = '  -10.53%  '
$ws.Range("D32").Value = '689.24'
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").Value = '11.81'
$ws.Range("E33").Value = '  -3.95%  '
$ws.Range("D34").Value = '64.27'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("E35").Value = '  -6.39%  '
$ws.Range("D36").Value = '39.18'
$ws.Range("E36").Value = '  -8.70%  '
$ws.Range("D37").Value = '0.412'
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").Value = '3.04'
$ws.Range("E39").Value = '  -2.92%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.131'
$ws.Range("E40").Value = '  -2.35%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '3.082.43'
$ws.Range("E41").Value = '  -3.16%  '
$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  -0.18%  '
$ws.Range("D43").Value = '0.0₃0700'
$ws.Range("E43").Value = '  -12.18%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.57'
$ws.Range("E44").Value = '  -12.60%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.83'
$ws.Range("E45").Value = '  +10.57%  '
$ws.Range("E46").Value = '  +6.84%  '
$ws.Range("D47").Value = '0.0400'
$ws.Range("E47").Value = '  -5.20%  '
$ws.Range("D48").Value = '0.127'
$ws.Range("E48").Value = '  -3.93%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '3.00'
$ws.Range("E49").Value = '  -3.83%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '137.25'
$ws.Range("E50").Value = '  -4.14%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '8.32'
$ws.Range("E51").Value = '  -6.13%  '
